$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (8) with the data for the 7th question, mirroring the
# existing Id / Question / Answer rows already in the sheet.
$ws.Range("A8").Value = "2023-10-22 21:36:58 7 question_7_8530875"
$ws.Range("B8").Value = "7 question"

# "7" looks numeric, so force the cell to text first (matching the other
# rows, e.g. C2="1", C3="2", ... which are stored as text) and then drop
# back to the workbook's default "Normal" style so no extra formatting is
# left behind on the cell.
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "7"
$ws.Range("C8").Style = "Normal"
